{"js": "// Append a new paragraph after the last existing paragraph in the body.\n// The new paragraph contains three separate (plain, unformatted) runs:\n//   \"examples of writing changes\", \"this text is bold\", \"this text is italic\"\n// We build the new paragraph with raw OOXML (via Range.insertOoxml, the\n// real Word JS API for inserting canonical WordprocessingML) so the\n// paragraph does NOT inherit the preceding paragraph's \"Title\" style and so\n// the three sentences land in three distinct <w:r> runs instead of being\n// merged into one run by repeated insertText() calls.\nconst body = context.document.body;\nbody.load(\"paragraphs\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// \"Flat OPC\" wrapper required by Range.insertOoxml.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>examples of writing changes</w:t></w:r>\n            <w:r><w:t>this text is bold</w:t></w:r>\n            <w:r><w:t>this text is italic</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst insertionRange = lastParagraph.getRange(\"After\");\ninsertionRange.insertOoxml(ooxml, \"After\");\nawait context.sync();\n", "ps1": "# Append a new paragraph after the last existing paragraph in the document.\n# The new paragraph holds three separate (plain, unformatted) runs:\n#   \"examples of writing changes\", \"this text is bold\", \"this text is italic\"\n#\n# We insert a new empty paragraph first (via Range.InsertParagraphAfter, so\n# it does not disturb the existing \"fisrt heading\" paragraph), then write the\n# three sentences into that new paragraph's Range as raw WordprocessingML via\n# Range.InsertXML. Using InsertXML (rather than three Range.InsertAfter /\n# TypeText calls) guarantees three distinct <w:r> runs and keeps the new\n# paragraph free of the preceding paragraph's \"Title\" style / <w:pPr>.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$tailRange = $lastParagraph.Range\n$tailRange.Collapse(0)          # wdCollapseEnd\n$tailRange.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$newRange = $newParagraph.Range\n$newParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>examples of writing changes</w:t></w:r><w:r><w:t>this text is bold</w:t></w:r><w:r><w:t>this text is italic</w:t></w:r></w:p>'\n$newRange.InsertXML($newParagraphXml)\n"}
